$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from Veneninvasion Co" sheet to "Include #0"
$sheet2 = $wb.Worksheets.Item("Include from Veneninvasion Co")
$sheet2.Name = "Include #0"

$ws = $wb.Worksheets.Item("Metadata")

# 2. Update the "Date" metadata value
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# 3. Insert a new "Jurisdiction" row right after "Contact" (row 10), which
#    pushes Description/Purpose/Copyright/Immutable down by one row each.
#    Capture the current row values first (read with the Value() call form).
$a11 = $ws.Range("A11").Value()
$b11 = $ws.Range("B11").Value()
$a12 = $ws.Range("A12").Value()
$b12 = $ws.Range("B12").Value()
$a13 = $ws.Range("A13").Value()
$b13 = $ws.Range("B13").Value()
$a14 = $ws.Range("A14").Value()
$b14 = $ws.Range("B14").Value()

# Row 15 is brand new on this sheet - clone row 14's formatting onto it
# before writing values so it keeps the same cell style.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

# Shift the captured rows down by one (bottom-up so nothing is clobbered).
$ws.Range("A15").Value = $a14
$ws.Range("B15").Value = $b14

$ws.Range("A14").Value = $a13
$ws.Range("B14").Value = $b13

$ws.Range("A13").Value = $a12
$ws.Range("B13").Value = $b12

$ws.Range("A12").Value = $a11
$ws.Range("B12").Value = $b11

# Write the new row's content. (The Jurisdiction value is an empty string
# in the source data - writing "" clears the cell, which is the closest
# representation reachable through the Value setter.)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
